# Fruta / hortaliza, semanal
#
# Insert a new weekly price-report row for "Uva" (Red Globe, Provincia de
# Limarí) at row 142 of the only worksheet, shifting the existing rows
# 142:227 down to 143:228 (the former last row, 227, becomes row 228).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142; existing rows 142:227 shift down to 143:228.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new data entry.
$ws.Cells.Item(142, 1).Value = 4
$ws.Cells.Item(142, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(142, 3).Value = "Los Lagos"
$ws.Cells.Item(142, 4).Value = 44680
$ws.Cells.Item(142, 5).Value = 10
$ws.Cells.Item(142, 6).Value = "Fruta"
$ws.Cells.Item(142, 7).Value = 100109
$ws.Cells.Item(142, 8).Value = "Uva"
$ws.Cells.Item(142, 9).Value = 100109001
$ws.Cells.Item(142, 10).Value = "Uva"
$ws.Cells.Item(142, 11).Value = "Red Globe"
$ws.Cells.Item(142, 12).Value = "Primera"
$ws.Cells.Item(142, 13).Value = 300
$ws.Cells.Item(142, 14).Value = 14000
$ws.Cells.Item(142, 15).Value = 15000
$ws.Cells.Item(142, 16).Value = 14500
$ws.Cells.Item(142, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(142, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(142, 19).Value = 725
$ws.Cells.Item(142, 20).Value = 20
